# For every color-named worksheet, fill column B (rows 2-15, the "Values"
# column) with the worksheet's own name, since each sheet represents one
# color and the cells were left blank placeholders in the source file.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $colorName = $ws.Name
    for ($row = 2; $row -le 15; $row++) {
        $ws.Cells.Item($row, 2).Value = $colorName
    }
}
